$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fill in the new bitacora entry on row 23 (C23 topic, D23 date, E23 comment).
# C23/E23 already carry the correct "label" style from the blank template row,
# so a plain value assignment is enough. D23 needs the date number-format that
# the other date cells (e.g. D21) use, so copy that formatting over first.
$ws.Range("C23").Value = "Modelos"

$ws.Range("D21").Copy()
$ws.Range("D23").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D23").Value = 46931

$ws.Range("E23").Value = "Creacion de modelo de modelos expenses y categorias"

# Update the view: scroll the window down a bit and move the active selection.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I22").Select()
